# configuración dns - dominio mileto.com.co
# Refresh the timestamp columns on the three demand-forecast sheets.
# Every value advances by the same 40 days, 1 hour, 18 minutes offset
# that the upstream data-refresh job applied.

$wb = $excel.ActiveWorkbook

# --- Sheet "24_HRS": column A holds bare HH:MM times (rows 2-25) ---
$ws24 = $wb.Worksheets.Item("24_HRS")
$hrs24 = @{
    2  = "21:40"
    3  = "20:40"
    4  = "19:40"
    5  = "18:40"
    6  = "17:40"
    7  = "16:40"
    8  = "15:40"
    9  = "14:40"
    10 = "13:40"
    11 = "12:40"
    12 = "11:40"
    13 = "10:40"
    14 = "09:40"
    15 = "08:40"
    16 = "07:40"
    17 = "06:40"
    18 = "05:40"
    19 = "04:40"
    20 = "03:40"
    21 = "02:40"
    22 = "01:40"
    23 = "00:40"
    24 = "23:40"
    25 = "22:40"
}
foreach ($row in $hrs24.Keys) {
    $ws24.Cells.Item($row, 1).Value = $hrs24[$row]
}

# --- Sheet "1d_bef": column A holds a single "YYYY-MM-DD HH:MM" timestamp ---
$ws1d = $wb.Worksheets.Item("1d_bef")
$ws1d.Cells.Item(2, 1).Value = "2023-05-02 22:40"

# --- Sheet "7d_bef": column A holds "YYYY-MM-DD HH:MM" timestamps (rows 2-25) ---
$ws7d = $wb.Worksheets.Item("7d_bef")
$dates7d = @{
    2  = "2023-04-26 22:40"
    3  = "2023-04-26 21:40"
    4  = "2023-04-26 20:40"
    5  = "2023-04-26 19:40"
    6  = "2023-04-26 18:40"
    7  = "2023-04-26 17:40"
    8  = "2023-04-26 16:40"
    9  = "2023-04-26 15:40"
    10 = "2023-04-26 14:40"
    11 = "2023-04-26 13:40"
    12 = "2023-04-26 12:40"
    13 = "2023-04-26 11:40"
    14 = "2023-04-26 10:40"
    15 = "2023-04-26 09:40"
    16 = "2023-04-26 08:40"
    17 = "2023-04-26 07:40"
    18 = "2023-04-26 06:40"
    19 = "2023-04-26 05:40"
    20 = "2023-04-26 04:40"
    21 = "2023-04-26 03:40"
    22 = "2023-04-26 02:40"
    23 = "2023-04-26 01:40"
    24 = "2023-04-26 00:40"
    25 = "2023-04-25 23:40"
}
foreach ($row in $dates7d.Keys) {
    $ws7d.Cells.Item($row, 1).Value = $dates7d[$row]
}
